# Auto-generated Excel COM-interop script
# Updates market-price derived columns (H:N) across several sheets
# to reflect a refreshed data pull from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1360
$ws.Range("I18").Value = 1420
$ws.Range("J18").Value = 1300
$ws.Range("K18").Value = 1420
$ws.Range("L18").Value = 1300
$ws.Range("M18").Value = -1136
$ws.Range("N18").Value = -1868

$ws.Range("H33").Value = 176.66667
$ws.Range("I33").Value = 159.2
$ws.Range("J33").Value = 198.5
$ws.Range("K33").Value = 159.2
$ws.Range("L33").Value = 198.5
$ws.Range("M33").Value = 69.80000000000001
$ws.Range("N33").Value = -656.5

$ws.Range("H53").Value = 367.33334
$ws.Range("I53").Value = 709.5
$ws.Range("J53").Value = 196.25
$ws.Range("K53").Value = 709.5
$ws.Range("L53").Value = 196.25
$ws.Range("M53").Value = -72.5
$ws.Range("N53").Value = -1470.25

$ws.Range("H100").Value = 1649.9
$ws.Range("I100").Value = 1749.8889
$ws.Range("J100").Value = 750
$ws.Range("K100").Value = 1749.8889
$ws.Range("L100").Value = 750
$ws.Range("M100").Value = -1208.8889
$ws.Range("N100").Value = -1832

$ws.Range("H113").Value = 5623.5
$ws.Range("I113").Value = 6113.8
$ws.Range("J113").Value = 4806.3335
$ws.Range("K113").Value = 6113.8
$ws.Range("L113").Value = 4806.3335
$ws.Range("M113").Value = -2859.8
$ws.Range("N113").Value = -11314.3335

$ws.Range("H132").Value = 12948.956
$ws.Range("I132").Value = 11469.789
$ws.Range("J132").Value = 19975
$ws.Range("K132").Value = 34409.367
$ws.Range("L132").Value = 59925
$ws.Range("M132").Value = -31879.367
$ws.Range("N132").Value = -64985

$ws.Range("H137").Value = 1778.6666
$ws.Range("I137").Value = 1294.25
$ws.Range("J137").Value = 2747.5
$ws.Range("K137").Value = 3882.75
$ws.Range("L137").Value = 8242.5
$ws.Range("M137").Value = -1332.75
$ws.Range("N137").Value = -13342.5

$ws.Range("H138").Value = 3267
$ws.Range("I138").Value = 914.7778
$ws.Range("J138").Value = 4895.4614
$ws.Range("K138").Value = 2744.3334
$ws.Range("L138").Value = 14686.3842
$ws.Range("M138").Value = 2395.6666
$ws.Range("N138").Value = -24966.3842

$ws.Range("H141").Value = 3065.5557
$ws.Range("I141").Value = 2416.6667
$ws.Range("J141").Value = 4363.3335
$ws.Range("K141").Value = 7250.000100000001
$ws.Range("L141").Value = 13090.0005
$ws.Range("M141").Value = -2070.000100000001
$ws.Range("N141").Value = -23450.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 1538.6666
$ws.Range("I24").Value = 1408
$ws.Range("J24").Value = 1800
$ws.Range("K24").Value = 1408
$ws.Range("L24").Value = 1800
$ws.Range("M24").Value = -1173
$ws.Range("N24").Value = -2270

$ws.Range("H64").Value = 709
$ws.Range("I64").Value = 777
$ws.Range("J64").Value = 675
$ws.Range("K64").Value = 777
$ws.Range("L64").Value = 675
$ws.Range("M64").Value = -552
$ws.Range("N64").Value = -1125

$ws.Range("H67").Value = 709
$ws.Range("I67").Value = 777
$ws.Range("J67").Value = 675
$ws.Range("K67").Value = 777
$ws.Range("L67").Value = 675
$ws.Range("M67").Value = 3
$ws.Range("N67").Value = -2235

$ws.Range("H86").Value = 7272.1816
$ws.Range("I86").Value = 5547.5
$ws.Range("J86").Value = 7655.4443
$ws.Range("K86").Value = 5547.5
$ws.Range("L86").Value = 7655.4443
$ws.Range("M86").Value = -4424.5
$ws.Range("N86").Value = -9901.444299999999

$ws.Range("H89").Value = 7272.1816
$ws.Range("I89").Value = 5547.5
$ws.Range("J89").Value = 7655.4443
$ws.Range("K89").Value = 27737.5
$ws.Range("L89").Value = 38277.2215
$ws.Range("M89").Value = -22121.5
$ws.Range("N89").Value = -49509.2215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = $null

$ws.Range("H31").Value = 4426.206
$ws.Range("I31").Value = 1957.4783
$ws.Range("J31").Value = 9588.091
$ws.Range("K31").Value = 1957.4783
$ws.Range("L31").Value = 9588.091
$ws.Range("M31").Value = -1662.4783
$ws.Range("N31").Value = -10178.091

$ws.Range("H34").Value = 4426.206
$ws.Range("I34").Value = 1957.4783
$ws.Range("J34").Value = 9588.091
$ws.Range("K34").Value = 1957.4783
$ws.Range("L34").Value = 9588.091
$ws.Range("M34").Value = -1755.4783
$ws.Range("N34").Value = -9992.091

$ws.Range("H58").Value = 4894.3
$ws.Range("I58").Value = 4400.6
$ws.Range("J58").Value = 5388
$ws.Range("K58").Value = 4400.6
$ws.Range("L58").Value = 5388
$ws.Range("M58").Value = -4197.6
$ws.Range("N58").Value = -5794

$ws.Range("H99").Value = 2346.75
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2396.2856
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2396.2856
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -5392.2856

$ws.Range("H107").Value = 295.77777
$ws.Range("I107").Value = 207.14285
$ws.Range("J107").Value = 606
$ws.Range("K107").Value = 207.14285
$ws.Range("L107").Value = 606
$ws.Range("M107").Value = 1712.85715
$ws.Range("N107").Value = -4446

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = $null

$ws.Range("H126").Value = 2346.75
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2396.2856
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 7188.8568
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -12128.8568

$ws.Range("H136").Value = 4894.3
$ws.Range("I136").Value = 4400.6
$ws.Range("J136").Value = 5388
$ws.Range("K136").Value = 13201.8
$ws.Range("L136").Value = 16164
$ws.Range("M136").Value = -10651.8
$ws.Range("N136").Value = -21264

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 79.47619
$ws.Range("I2").Value = 99.28570999999999
$ws.Range("J2").Value = 39.857143
$ws.Range("K2").Value = 595.71426
$ws.Range("L2").Value = 239.142858
$ws.Range("M2").Value = -482.71426
$ws.Range("N2").Value = -465.142858

$ws.Range("H23").Value = 176.6
$ws.Range("I23").Value = 85.333336
$ws.Range("J23").Value = 215.71428
$ws.Range("K23").Value = 256.000008
$ws.Range("L23").Value = 647.14284
$ws.Range("M23").Value = -21.00000799999998
$ws.Range("N23").Value = -1117.14284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 370.4
$ws.Range("I9").Value = 370.4
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 370.4
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -200.4

$ws.Range("H70").Value = 3709.6
$ws.Range("I70").Value = 3665.375
$ws.Range("J70").Value = 3886.5
$ws.Range("K70").Value = 3665.375
$ws.Range("L70").Value = 3886.5
$ws.Range("M70").Value = -3395.375
$ws.Range("N70").Value = -4426.5

$ws.Range("H73").Value = 3709.6
$ws.Range("I73").Value = 3665.375
$ws.Range("J73").Value = 3886.5
$ws.Range("K73").Value = 3665.375
$ws.Range("L73").Value = 3886.5
$ws.Range("M73").Value = -2729.375
$ws.Range("N73").Value = -5758.5

$ws.Range("H102").Value = 1941.1111
$ws.Range("I102").Value = 1823.4615
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 1823.4615
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -201.4614999999999
$ws.Range("N102").Value = -8244

$ws.Range("H132").Value = 3096.6924
$ws.Range("I132").Value = 2521.4167
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 7564.250100000001
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -5034.250100000001
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 150
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 150
$ws.Range("N22").Value = -740

$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 150
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 150
$ws.Range("N27").Value = -364

$ws.Range("H46").Value = 6348.75
$ws.Range("I46").Value = 2600
$ws.Range("J46").Value = 8598
$ws.Range("K46").Value = 2600
$ws.Range("L46").Value = 8598
$ws.Range("M46").Value = -2412
$ws.Range("N46").Value = -8974

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = $null
$ws.Range("N96").Value = $null
